# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the per-market-day records (rows 2-18) of the
# "Haba" sheet: each row keeps its fixed attributes (Mercado, Region, Codreg,
# Categoria, Variedad, Calidad, Unidad de comercializacion, Kg o Unidades,
# Clasificacion) but is re-paired with a different day's figures (Fecha,
# Volumen, Precio minimo/maximo/promedio ponderado, Origen, Precio $/Kg).
# Row 16 is unaffected by this week's reshuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
#          M=Precio promedio ponderado, O=Origen, P=Precio $/Kg

$ws.Cells.Item(2, 4).Value = 44446
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 11000
$ws.Cells.Item(2, 13).Value = 11500
$ws.Cells.Item(2, 16).Value = 460

$ws.Cells.Item(3, 4).Value = 44316
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 16000
$ws.Cells.Item(3, 12).Value = 17000
$ws.Cells.Item(3, 13).Value = 16500
$ws.Cells.Item(3, 16).Value = 660

$ws.Cells.Item(4, 4).Value = 44370
$ws.Cells.Item(4, 10).Value = 520
$ws.Cells.Item(4, 11).Value = 13000
$ws.Cells.Item(4, 12).Value = 14000
$ws.Cells.Item(4, 13).Value = 13500
$ws.Cells.Item(4, 16).Value = 540

$ws.Cells.Item(5, 4).Value = 44356
$ws.Cells.Item(5, 10).Value = 500
$ws.Cells.Item(5, 15).Value = 'Provincia de Limarí'

$ws.Cells.Item(6, 4).Value = 44694
$ws.Cells.Item(6, 10).Value = 480
$ws.Cells.Item(6, 11).Value = 17500
$ws.Cells.Item(6, 12).Value = 18000
$ws.Cells.Item(6, 13).Value = 17750
$ws.Cells.Item(6, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(6, 16).Value = 710

$ws.Cells.Item(7, 4).Value = 44690
$ws.Cells.Item(7, 10).Value = 400
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17500
$ws.Cells.Item(7, 16).Value = 700

$ws.Cells.Item(8, 4).Value = 44377
$ws.Cells.Item(8, 10).Value = 520
$ws.Cells.Item(8, 11).Value = 12500
$ws.Cells.Item(8, 13).Value = 12750
$ws.Cells.Item(8, 16).Value = 510

$ws.Cells.Item(9, 4).Value = 44721
$ws.Cells.Item(9, 10).Value = 500
$ws.Cells.Item(9, 11).Value = 14500
$ws.Cells.Item(9, 12).Value = 15000
$ws.Cells.Item(9, 13).Value = 14750
$ws.Cells.Item(9, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(9, 16).Value = 590

$ws.Cells.Item(10, 4).Value = 44484
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 9000
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 13).Value = 9500
$ws.Cells.Item(10, 16).Value = 380

$ws.Cells.Item(11, 4).Value = 44386
$ws.Cells.Item(11, 10).Value = 500
$ws.Cells.Item(11, 11).Value = 11000
$ws.Cells.Item(11, 12).Value = 12000
$ws.Cells.Item(11, 13).Value = 11500
$ws.Cells.Item(11, 16).Value = 460

$ws.Cells.Item(12, 4).Value = 44372
$ws.Cells.Item(12, 10).Value = 500
$ws.Cells.Item(12, 11).Value = 13000
$ws.Cells.Item(12, 12).Value = 14000
$ws.Cells.Item(12, 13).Value = 13500
$ws.Cells.Item(12, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(12, 16).Value = 540

$ws.Cells.Item(13, 4).Value = 44384
$ws.Cells.Item(13, 10).Value = 560
$ws.Cells.Item(13, 11).Value = 11500
$ws.Cells.Item(13, 12).Value = 12000
$ws.Cells.Item(13, 13).Value = 11750
$ws.Cells.Item(13, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(13, 16).Value = 470

$ws.Cells.Item(14, 4).Value = 44425
$ws.Cells.Item(14, 10).Value = 400
$ws.Cells.Item(14, 11).Value = 11500
$ws.Cells.Item(14, 12).Value = 12000
$ws.Cells.Item(14, 13).Value = 11750
$ws.Cells.Item(14, 16).Value = 470

$ws.Cells.Item(15, 4).Value = 44466
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 9500
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 9750
$ws.Cells.Item(15, 16).Value = 390

$ws.Cells.Item(17, 4).Value = 44714
$ws.Cells.Item(17, 10).Value = 400
$ws.Cells.Item(17, 11).Value = 14000
$ws.Cells.Item(17, 12).Value = 15000
$ws.Cells.Item(17, 13).Value = 14500
$ws.Cells.Item(17, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(17, 16).Value = 580

$ws.Cells.Item(18, 4).Value = 44376
$ws.Cells.Item(18, 11).Value = 12000
$ws.Cells.Item(18, 12).Value = 13000
$ws.Cells.Item(18, 13).Value = 12500
$ws.Cells.Item(18, 16).Value = 500
